$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(34)
$src.Copy([System.Reflection.Missing]::Value, $src)
$ws = $wb.Worksheets.Item(35)
$ws.Name = '2025-07-26'

$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = '異世界おじさん'
$ws.Cells.Item(2,3).Value = '殆ど死んでいる(著者)'
$ws.Cells.Item(2,4).Value = '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = '悪人面したＢ級冒険者 主人公とその幼馴染たちのパパになる'
$ws.Cells.Item(3,3).Value = 'こげめ(著者) えんじ(原作) ハラカズヒロ(キャラクター原案)'
$ws.Cells.Item(3,4).Value = '「名もなき英雄譚」後半'
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = '淫獄団地'
$ws.Cells.Item(4,3).Value = '搾精研究所(原作) 丈山雄為(漫画)'
$ws.Cells.Item(4,4).Value = '第49話（前編）'
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = 'まんきつしたい常連さん'
$ws.Cells.Item(5,3).Value = 'しんみりん(著者)'
$ws.Cells.Item(5,4).Value = '第46話前編'
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = '男嫌いな美人姉妹を名前も告げずに助けたら一体どうなる?'
$ws.Cells.Item(6,3).Value = 'みょん(原作) 司馬淳子(漫画) ぎうにう(キャラクターデザイン)'
$ws.Cells.Item(6,4).Value = '第23話'
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = '美人女上司滝沢さん'
$ws.Cells.Item(7,3).Value = 'やんBARU(著者)'
$ws.Cells.Item(7,4).Value = '第202話'
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = 'よくわからないけれど異世界に転生していたようです'
$ws.Cells.Item(8,3).Value = '内々けやき あし カオミン'
$ws.Cells.Item(8,4).Value = '第136話 よくわからないけれどスカウトされたみたいです（２）'
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = '治癒魔法の間違った使い方 ~戦場を駆ける回復要員~'
$ws.Cells.Item(9,3).Value = '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)'
$ws.Cells.Item(9,4).Value = '第80話その3'
$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = '解雇された暗黒兵士(30代)のスローなセカンドライフ'
$ws.Cells.Item(10,3).Value = '岡沢六十四 るれくちぇ sage・ジョー'
$ws.Cells.Item(10,4).Value = '第71話(前編) ダリエルVS.滾り'
$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = '老後に備えて異世界で８万枚の金貨を貯めます'
$ws.Cells.Item(11,3).Value = 'FUNA 東西 モトエ恵介'
$ws.Cells.Item(11,4).Value = '第120話　会談［その5］'
$ws.Cells.Item(12,1).Value = 11
$ws.Cells.Item(12,2).Value = '天獄で悪魔がボクを魅惑する'
$ws.Cells.Item(12,3).Value = '銀河味めてお(著者)'
$ws.Cells.Item(12,4).Value = '第34話'
$ws.Cells.Item(13,1).Value = 12
$ws.Cells.Item(13,2).Value = 'ワンパンマン'
$ws.Cells.Item(13,3).Value = '原作/ＯＮＥ 作画/村田雄介'
$ws.Cells.Item(13,4).Value = '206撃目'
$ws.Cells.Item(14,1).Value = 13
$ws.Cells.Item(14,2).Value = 'ずっと好きだった幼馴染と付き合い始めたら一途ビッチの性欲ジャンキーだったんだがどうすりゃいいですか？'
$ws.Cells.Item(14,3).Value = '原作：トラ子猫 作画：あらいぐま'
$ws.Cells.Item(14,4).Value = '第3話'
$ws.Cells.Item(15,1).Value = 14
$ws.Cells.Item(15,2).Value = '冒険者絶対殺すダンジョン'
$ws.Cells.Item(15,3).Value = '道満晴明(著者)'
$ws.Cells.Item(15,4).Value = '第34話'
$ws.Cells.Item(16,1).Value = 15
$ws.Cells.Item(16,2).Value = 'ノロマ魔法と呼ばれた魔法使いは重力魔法で無双する　～まだ重力の概念のない世界にて、少年は万有引力の王となる～'
$ws.Cells.Item(16,3).Value = '神原絵理華(漫画) 一森一輝(原作)'
$ws.Cells.Item(16,4).Value = '第18話②'
$ws.Cells.Item(17,1).Value = 16
$ws.Cells.Item(17,2).Value = '両親の借金を肩代わりしてもらう条件は日本一可愛い女子高生と一緒に暮らすことでした。'
$ws.Cells.Item(17,3).Value = '美月めいあ(漫画) 雨音恵(原作) ｋａｋａｏ(キャラクター原案)'
$ws.Cells.Item(17,4).Value = '第36話'
$ws.Cells.Item(18,1).Value = 17
$ws.Cells.Item(18,2).Value = '最強で最速の無限レベルアップ ～スキル【経験値1000倍】と【レベルフリー】でレベル上限の枷が外れた俺は無双する～'
$ws.Cells.Item(18,3).Value = 'シオヤマ琴 鳥羽田 航 トモゼロ'
$ws.Cells.Item(18,4).Value = '第75話 ゆめうつつ'
$ws.Cells.Item(19,1).Value = 18
$ws.Cells.Item(19,2).Value = '陰キャの僕に罰ゲームで告白してきたはずのギャルが、どう見ても僕にベタ惚れです'
$ws.Cells.Item(19,3).Value = '神奈なごみ(漫画) 結石(原作) かがちさく(キャラクター原案)'
$ws.Cells.Item(19,4).Value = '第26話'
$ws.Cells.Item(20,1).Value = 19
$ws.Cells.Item(20,2).Value = '修羅幼女の英雄譚～半端者と言われた傭兵、幼女に転生して成り上がる～'
$ws.Cells.Item(20,3).Value = '作画：むらたん 原作：沙城流'
$ws.Cells.Item(20,4).Value = '第7話(3)'
$ws.Cells.Item(21,1).Value = 20
$ws.Cells.Item(21,2).Value = 'ゴリラ女子高生'
$ws.Cells.Item(21,3).Value = '大友しゅうま(著者)'
$ws.Cells.Item(21,4).Value = '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'
$ws.Cells.Item(22,1).Value = 21
$ws.Cells.Item(22,2).Value = '王子様の友達'
$ws.Cells.Item(22,3).Value = 'すけろく(著者)'
$ws.Cells.Item(22,4).Value = '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'
$ws.Cells.Item(23,1).Value = 22
$ws.Cells.Item(23,2).Value = 'ラブコメと怪獣退治の不文律'
$ws.Cells.Item(23,3).Value = '御池慧（漫画） 上代務（原作） TMSLab（原作）'
$ws.Cells.Item(23,4).Value = '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'
$ws.Cells.Item(24,1).Value = 23
$ws.Cells.Item(24,2).Value = 'わたし、二番目の彼女でいいから。'
$ws.Cells.Item(24,3).Value = 'にの子(漫画) 西条陽(原作) Re岳(キャラクター原案)'
$ws.Cells.Item(24,4).Value = '第10話③：カーテンのなか'
$ws.Cells.Item(25,1).Value = 24
$ws.Cells.Item(25,2).Value = 'ヤンデレかと思ったらもっとヤベー女だった'
$ws.Cells.Item(25,3).Value = '八木戸マト'
$ws.Cells.Item(25,4).Value = '第68話　全てを失ったヤンデレ彼女'
$ws.Cells.Item(26,1).Value = 25
$ws.Cells.Item(26,2).Value = 'センパイ、自宅警備員の雇用はいかがですか？'
$ws.Cells.Item(26,3).Value = '漫画：コブラサナギ 原作：二上圭 キャラ原案：日向あずり'
$ws.Cells.Item(26,4).Value = '第5話後半'
$ws.Cells.Item(27,1).Value = 26
$ws.Cells.Item(27,2).Value = 'スキルがなければレベルを上げる～９９がカンストの世界でレベル800万からスタート～'
$ws.Cells.Item(27,3).Value = '倉橋ユウス(漫画) 岡沢六十四(原作)'
$ws.Cells.Item(27,4).Value = '第51話④'
$ws.Cells.Item(28,1).Value = 27
$ws.Cells.Item(28,2).Value = '時森さんが無防備です!!'
$ws.Cells.Item(28,3).Value = 'たざわ'
$ws.Cells.Item(28,4).Value = '第63話'
$ws.Cells.Item(29,1).Value = 28
$ws.Cells.Item(29,2).Value = '不老不死少女の苗床旅行記'
$ws.Cells.Item(29,3).Value = 'ふじはん(漫画) ルナ・ウサギ(原作)'
$ws.Cells.Item(29,4).Value = '第16話後編'
$ws.Cells.Item(30,1).Value = 29
$ws.Cells.Item(30,2).Value = 'くじ引き特賞：無双ハーレム権'
$ws.Cells.Item(30,3).Value = '原作／三木なずな（GA文庫／SBクリエイティブ刊） 漫画／長谷見亮 キャラクター原案／瑠奈璃亜'
$ws.Cells.Item(30,4).Value = '第58話-02　新たな王女たちへ、受け継がれし慈愛の心！'
$ws.Cells.Item(31,1).Value = 30
$ws.Cells.Item(31,2).Value = '勇者のクズ'
$ws.Cells.Item(31,3).Value = 'ナカシマ723'
$ws.Cells.Item(31,4).Value = '第48話　勇者の矜持 II（前半）'
$ws.Cells.Item(32,1).Value = 31
$ws.Cells.Item(32,2).Value = '婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版'
$ws.Cells.Item(32,3).Value = '漫画/すたひろ 原作/Y.A'
$ws.Cells.Item(32,4).Value = 'chapter66【35話①】'
$ws.Cells.Item(33,1).Value = 32
$ws.Cells.Item(33,2).Value = '断れない会長は友江くんにだけしてあげたい'
$ws.Cells.Item(33,3).Value = '沼地どろまる(著者)'
$ws.Cells.Item(33,4).Value = '休載漫画'
$ws.Cells.Item(34,1).Value = 33
$ws.Cells.Item(34,2).Value = 'ダウナーお姉さんは遊びたい'
$ws.Cells.Item(34,3).Value = '山鷹景'
$ws.Cells.Item(34,4).Value = '第15話'
$ws.Cells.Item(35,1).Value = 34
$ws.Cells.Item(35,2).Value = '聖騎士になったけど団長のおっぱいが凄すぎて心が清められない'
$ws.Cells.Item(35,3).Value = '木の芽(原作) 川喜田ミツオ(漫画) 雨傘ゆん(キャラクター原案)'
$ws.Cells.Item(35,4).Value = 'コミックス6巻発売＆休載のお知らせ'
$ws.Cells.Item(36,1).Value = 35
$ws.Cells.Item(36,2).Value = '脱稿するまでオチません'
$ws.Cells.Item(36,3).Value = 'ヨシラギ(著者)'
$ws.Cells.Item(36,4).Value = '第32話後半'
$ws.Cells.Item(37,1).Value = 36
$ws.Cells.Item(37,2).Value = '勇者パーティから追い出された不遇職【罠士】、ユニークスキル【矢印】で最強になる'
$ws.Cells.Item(37,3).Value = '作画：たつひこ 原作：白石 有希'
$ws.Cells.Item(37,4).Value = '第7話(3)'
$ws.Cells.Item(38,1).Value = 37
$ws.Cells.Item(38,2).Value = 'クラスの大嫌いな女子と結婚することになった。'
$ws.Cells.Item(38,3).Value = '天乃聖樹(原作) もすこんぶ(漫画)'
$ws.Cells.Item(38,4).Value = '第44話-2'
$ws.Cells.Item(39,1).Value = 38
$ws.Cells.Item(39,2).Value = '隣のクーデレラを甘やかしたら、ウチの合鍵を渡すことになった'
$ws.Cells.Item(39,3).Value = '青島かなえ(漫画) 雪仁(原作) かがちさく(キャラクター原案)'
$ws.Cells.Item(39,4).Value = '第38話'
$ws.Cells.Item(40,1).Value = 39
$ws.Cells.Item(40,2).Value = 'りんちゃんは据え膳したい'
$ws.Cells.Item(40,3).Value = '澄田佑貴(著者)'
$ws.Cells.Item(40,4).Value = '第38話'
$ws.Cells.Item(41,1).Value = 40
$ws.Cells.Item(41,2).Value = 'ボクの理想の異世界生活 ～転生したらケモ耳娘だらけの世界でハーレムに～'
$ws.Cells.Item(41,3).Value = 'イチリ(原作) 空維深夜(作画)'
$ws.Cells.Item(41,4).Value = '第15話前半：CALL'
$ws.Cells.Item(42,1).Value = 41
$ws.Cells.Item(42,2).Value = 'その冒険者、取り扱い注意。 ～正体は無敵の下僕たちを統べる異世界最強の魔導王～'
$ws.Cells.Item(42,3).Value = '満月シオン(作画) Sin Guilty(ツギクル)(原作) M.B(キャラクター原案)'
$ws.Cells.Item(42,4).Value = '56章　はじまりの愚か者②　後編'
$ws.Cells.Item(43,1).Value = 42
$ws.Cells.Item(43,2).Value = '俺の愛娘は悪役令嬢'
$ws.Cells.Item(43,3).Value = 'かわもり かぐら(原作) ほづみりや(漫画) 縞(キャラクター原案)'
$ws.Cells.Item(43,4).Value = '第4話-2'
$ws.Cells.Item(44,1).Value = 43
$ws.Cells.Item(44,2).Value = 'フルメタル・パニック！　Family'
$ws.Cells.Item(44,3).Value = '賀東招二(原作) 神反ヲ鬚(作画) 四季童子(キャラクター原案)'
$ws.Cells.Item(44,4).Value = '第6話　東京都江東区のタワマン39階②-2'
$ws.Cells.Item(45,1).Value = 44
$ws.Cells.Item(45,2).Value = '双子の姉が神子として引き取られて、私は捨てられたけど多分私が神子である。'
$ws.Cells.Item(45,3).Value = '雪(著者) 池中織奈(原作) カット(キャラクター原案)'
$ws.Cells.Item(45,4).Value = '第32話後編'
$ws.Cells.Item(46,1).Value = 45
$ws.Cells.Item(46,2).Value = 'TRPGプレイヤーが異世界で最強ビルドを目指す ～ヘンダーソン氏の福音を～'
$ws.Cells.Item(46,3).Value = '内田テモ(漫画) Schuld(原作) ランサネ(キャラクター原案)'
$ws.Cells.Item(46,4).Value = '第14話'
$ws.Cells.Item(47,1).Value = 46
$ws.Cells.Item(47,2).Value = 'クロの戦記Ⅱ 異世界転移した僕が最強なのはベッドの上だけのようです'
$ws.Cells.Item(47,3).Value = 'サイトウアユム(原作) ユリシロ(漫画) むつみまさと(キャラクター原案)'
$ws.Cells.Item(47,4).Value = '第22話-2'
$ws.Cells.Item(48,1).Value = 47
$ws.Cells.Item(48,2).Value = 'ダークサモナーとデキている'
$ws.Cells.Item(48,3).Value = '車王(著者)'
$ws.Cells.Item(48,4).Value = '【コミックス第6巻発売記念】挟まるならどの衣装？コメント大募集！'
$ws.Cells.Item(49,1).Value = 48
$ws.Cells.Item(49,2).Value = 'ギルドを追放された回復術士、実は魔力無限だったので規格外の回復魔法で伝説となる'
$ws.Cells.Item(49,3).Value = '漫画：坂下コウ 原作：霞杏檎'
$ws.Cells.Item(49,4).Value = '第4話(2)'
$ws.Cells.Item(50,1).Value = 49
$ws.Cells.Item(50,2).Value = 'チンチンデビルを追え！'
$ws.Cells.Item(50,3).Value = 'くぼたふみお'
$ws.Cells.Item(50,4).Value = '第３１話　激突！ アンとキュベレ！'
$ws.Cells.Item(51,1).Value = 50
$ws.Cells.Item(51,2).Value = '俺の『鑑定』スキルがチートすぎて ～伝説の勇者を読み“盗り”最強へ～'
$ws.Cells.Item(51,3).Value = '原作：澄守　彩 漫画：龍牙 翔'
$ws.Cells.Item(51,4).Value = '第33話ー③　迷宮の掃除番'
